$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (72) down into the new
# row (73), so the new row picks up the same cell styles already used by
# the table (date format in column A, centered numbers in B:F) instead of
# creating brand-new style entries.
$ws.Range("A72:F72").Copy()
$ws.Range("A73:F73").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new day's figures (data published 05-25-2020 for 05-24-2020).
$ws.Range("A73").Value = 43975
$ws.Range("B73").Value = 667
$ws.Range("C73").Value = 207
$ws.Range("D73").Value = 381
$ws.Range("E73").Value = 26
$ws.Range("F73").Value = 20

# Grow the table / AutoFilter range to include the new row.
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:F73"))

# Mirror the author's final selection on the sheet.
$ws.Range("F73").Select()
